$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename existing columns ---
$ws.Range("B1").Value = "ID_Cliente"
$ws.Range("C1").Value = "ID_Servico"
$ws.Range("E1").Value = "ID_Funcionario"
$ws.Range("F1").Value = "Valor_Total"

# --- Header row (row 1): two new columns ---
$ws.Range("G1").Value = "Data_Ultimo_Atendimento"
$ws.Range("H1").Value = "Execucao"

# Give the two new header cells the same look as the existing header
# cells (bold font, centered/top alignment, thin border) by copying the
# formatting from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data row (row 2): update existing values ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = 1

# D2 holds a date written as plain text ("2025-02-06"), not a real date
# serial. Force the cell to Text format before assigning the string so
# Excel doesn't auto-convert it to a date serial number, then clear the
# formatting back off so the cell keeps the workbook's default (no
# explicit) style, matching the rest of the row.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2025-02-06"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 50

# --- Data row (row 2): two new values ---
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2025-02-06"
$ws.Range("G2").ClearFormats()

$ws.Range("H2").Value = "Andamento"
